# Update shop data workbook:
#  - Products, Sales, Expenses sheets get their header columns reordered
#    (and Expenses drops the "receipt_number" column)
#  - Each of those three sheets gains a new data row (row 2) representing
#    a product, the sale of that product, and an unrelated expense entry.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Products
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Products")

$headers1 = @("name","description","price","cost_price","category","stock","min_stock","supplier","sku","id","created_date","last_updated")
for ($i = 0; $i -lt $headers1.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).NumberFormat = "@"
    $ws.Cells.Item(1, $i + 1).Value = $headers1[$i]
}

$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "Amul Butter (500g)"

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = ""

$ws.Range("C2").NumberFormat = "General"
$ws.Range("C2").Value = 500

$ws.Range("D2").NumberFormat = "General"
$ws.Range("D2").Value = 100

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "Food"

$ws.Range("F2").NumberFormat = "General"
$ws.Range("F2").Value = 29

$ws.Range("G2").NumberFormat = "General"
$ws.Range("G2").Value = 1

$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = ""

$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = ""

$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "07291cea-a90f-4352-b915-1121e63dcb59"

$ws.Range("K2").NumberFormat = "@"
$ws.Range("K2").Value = "2025-09-23T12:50:58.115Z"

$ws.Range("L2").NumberFormat = "@"
$ws.Range("L2").Value = "2025-09-23T12:51:14.513Z"

# ---------------------------------------------------------------------
# Sheet 2: Sales
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Sales")

$headers2 = @("product_id","quantity","unit_price","customer_name","payment_method","cashier","notes","id","product_name","total_amount","profit","sale_date")
for ($i = 0; $i -lt $headers2.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).NumberFormat = "@"
    $ws.Cells.Item(1, $i + 1).Value = $headers2[$i]
}

$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "07291cea-a90f-4352-b915-1121e63dcb59"

$ws.Range("B2").NumberFormat = "General"
$ws.Range("B2").Value = 71

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "500"

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = ""

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "Cash"

$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "Admin"

$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = ""

$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "9cbf4f85-2a79-4a90-8951-8c7ed4927e8d"

$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "Amul Butter (500g)"

$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "35500"

$ws.Range("K2").NumberFormat = "@"
$ws.Range("K2").Value = "28400"

$ws.Range("L2").NumberFormat = "@"
$ws.Range("L2").Value = "2025-09-23T12:51:14.506Z"

# ---------------------------------------------------------------------
# Sheet 3: Expenses
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Expenses")

$headers3 = @("category","description","amount","payment_method","vendor","notes","id","expense_date")
for ($i = 0; $i -lt $headers3.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).NumberFormat = "@"
    $ws.Cells.Item(1, $i + 1).Value = $headers3[$i]
}

# Remove the stale "receipt_number" column data if it spilled past the new
# header width (old sheet had 9 columns, new one has 8).
$ws.Range("I1").Clear()

$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "Utilities"

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "Random"

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "5000"

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "Cash"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = ""

$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = ""

$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "0f1d1996-63c0-461d-b044-35309c9cda10"

$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "2025-09-23T12:51:31.350Z"
